$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.133.07'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '2.049.34'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  -0.03%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.14'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -2.29%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.663'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -2.04%  '
$ws.Range("E7").Value = '  +0.04%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.22'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -7.43%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.381'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  -2.87%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0782'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -2.81%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.109'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -0.31%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.33'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +0.44%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.881'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +7.41%  '
$ws.Range("D14").Value = '2.346.43'
$ws.Range("E14").Value = '  -1.44%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.72'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("D16").Value = '2.051.60'
$ws.Range("E16").Value = '  -1.46%  '
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.34'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +10.01%  '
$ws.Range("D18").Value = '37.139.87'
$ws.Range("E18").Value = '  -0.11%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.54'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").Value = '0.0₃0894'
$ws.Range("E20").Value = '  -4.50%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.41'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -1.91%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '236.87'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("E23").Value = '  -0.02%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.48'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +2.43%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.55'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +1.13%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.73'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -0.31%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.17'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -5.39%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.10'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -1.79%  '
$ws.Range("E29").Value = '  -1.81%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.92'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +2.17%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  -0.35%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0619'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -1.37%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.48'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("E35").Value = '  -0.05%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.25'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -3.16%  '
$ws.Range("E37").Value = '  +0.57%  '
$ws.Range("E38").Value = '  -3.12%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.26'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +13.98%  '
$ws.Range("E40").Value = '  +9.39%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0986'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -17.49%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0223'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -2.53%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.45'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("E44").Value = '  -2.02%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '95.68'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -3.66%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.43'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -3.21%  '
$ws.Range("D47").Value = '1.266.40'
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("E48").Value = '  -3.00%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.77'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -2.99%  '
$ws.Range("D50").Value = '2.229.68'
$ws.Range("E50").Value = '  -1.62%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.04'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -1.35%  '
